# "Insert row" action for the transmittal sheet: appends a new data row
# right after the current last row (row 15 -> new row 16), matching the
# sheet's existing look-and-feel, and grows the used range accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 15
$newRow  = $lastRow + 1

# Duplicate the (currently blank/unformatted) last row down into the new
# row first. Columns C and F carry a column-level default style in this
# sheet, and simply typing a value straight into a brand-new cell in those
# columns would silently pick that style up - pre-creating the row-16
# cells as plain copies of row 15 keeps them unstyled, same as the rest of
# the newly-added row.
$ws.Range("A" + $lastRow + ":F" + $lastRow).Copy()
$ws.Range("A" + $newRow + ":F" + $newRow).PasteSpecial(-4122)

# The previously-last row now becomes a normal interior row, so give it
# the same formatting already used by the rows above it (13/14).
$ws.Rows($lastRow).Style = $ws.Rows($lastRow - 2).Style

# Fill in the newly-inserted row's values (Control #, DV, Check Date,
# Check #, Payee, Status) - the last "Amount" column is left blank, as
# is the Status column's default for this particular entry.
$ws.Cells.Item($newRow, 1).Value = "asdf"
$ws.Cells.Item($newRow, 2).Value = "asdf"
$ws.Cells.Item($newRow, 3).Value = "sadf"
$ws.Cells.Item($newRow, 4).Value = "asdfasdf"
$ws.Cells.Item($newRow, 5).Value = "asdf"
$ws.Cells.Item($newRow, 6).Value = "Paid"
